# Apply the "Add data for 2021-12-29" update:
#  - Rename the worksheet from "Through 2021-12-20" to "Through 2021-12-21"
#  - Update the label in A14 to match
#  - Update the December (row 14) and Total (row 15) figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Through 2021-12-21"

# Row 14 - "December (through 12-20)" -> "December (through 12-21)" and new figures
$ws.Range("A14").Value = "December (through 12-21)"
$ws.Range("C14").Value = 25
$ws.Range("D14").Value = 0.1071
$ws.Range("F14").Value = 60
$ws.Range("G14").Value = 0.09089999999999999
$ws.Range("I14").Value = 72
$ws.Range("J14").Value = 0.1111
$ws.Range("L14").Value = 43
$ws.Range("M14").Value = 0.0851
$ws.Range("N14").Value = 4
$ws.Range("O14").Value = 35
$ws.Range("P14").Value = 0.1026
$ws.Range("Q14").Value = 5
$ws.Range("R14").Value = 93
$ws.Range("S14").Value = 0.051
$ws.Range("U14").Value = 137
$ws.Range("V14").Value = 0.0144

# Row 15 - Total figures
$ws.Range("C15").Value = 283
$ws.Range("D15").Value = 0.1129
$ws.Range("F15").Value = 564
$ws.Range("G15").Value = 0.1033
$ws.Range("I15").Value = 830
$ws.Range("J15").Value = 0.0798
$ws.Range("L15").Value = 651
$ws.Range("M15").Value = 0.107
$ws.Range("N15").Value = 58
$ws.Range("O15").Value = 515
$ws.Range("P15").Value = 0.1012
$ws.Range("Q15").Value = 69
$ws.Range("R15").Value = 1293
$ws.Range("S15").Value = 0.0507
$ws.Range("U15").Value = 1679
